$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9999999945686564
$ws.Range("A2").Value = 0.99817057862474745
$ws.Range("A3").Value = 0.99197263649637768
$ws.Range("A4").Value = 0.99327596452985389
$ws.Range("A5").Value = 0.98137576892656098
$ws.Range("A6").Value = 0.95251692698866008
$ws.Range("A7").Value = 0.94984072761529992
$ws.Range("A8").Value = 0.94637472550573865
$ws.Range("A9").Value = 0.94817779625088539
$ws.Range("A10").Value = 0.95096026816593371
$ws.Range("A11").Value = 0.95152081752405993
$ws.Range("A12").Value = 0.95110454821094859
$ws.Range("A13").Value = 0.95041237075784368
$ws.Range("A14").Value = 0.95071501791483259
$ws.Range("A15").Value = 0.94812392405350354
$ws.Range("A16").Value = 0.94561770765467679
$ws.Range("A17").Value = 0.94191004832525038
$ws.Range("A18").Value = 0.94080114554307048
$ws.Range("A19").Value = 0.99711736930638351
$ws.Range("A20").Value = 0.99000035946713871
$ws.Range("A21").Value = 0.98860187159373292
$ws.Range("A22").Value = 0.98733736851654208
$ws.Range("A23").Value = 0.98493834444107442
$ws.Range("A24").Value = 0.97191798516615324
$ws.Range("A25").Value = 0.96546107352209121
$ws.Range("A26").Value = 0.94659107841376522
$ws.Range("A27").Value = 0.94174561962528514
$ws.Range("A28").Value = 0.92027635848831457
$ws.Range("A29").Value = 0.90500625560284687
$ws.Range("A30").Value = 0.89843609747833897
$ws.Range("A31").Value = 0.8907824422276488
$ws.Range("A32").Value = 0.88910311290450816
$ws.Range("A33").Value = 0.88858310167505161
